$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.833.73'
$ws.Range("E2").Value = '  +1.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.934.98'
$ws.Range("E3").Value = '  +1.24%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.36'
$ws.Range("E5").Value = '  +4.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.008'
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4837'
$ws.Range("E7").Value = '  +0.37%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4109'
$ws.Range("E8").Value = '  +1.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08176'
$ws.Range("E9").Value = '  -0.54%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.014'
$ws.Range("E10").Value = '  -0.69%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.67'
$ws.Range("E11").Value = '  +0.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.986.39'
$ws.Range("E12").Value = '  +3.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.079'
$ws.Range("E13").Value = '  +0.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.274'
$ws.Range("E14").Value = '  +1.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.89'
$ws.Range("E15").Value = '  -0.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06844'
$ws.Range("E16").Value = '  +0.64%  '

$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001034'
$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.78'
$ws.Range("E19").Value = '  +0.38%  '

$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.828.61'
$ws.Range("E21").Value = '  +1.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.631'
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.90'
$ws.Range("E23").Value = '  +0.85%  '

$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.235.43'
$ws.Range("E24").Value = '  +4.37%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.181'
$ws.Range("E25").Value = '  -0.22%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.10'
$ws.Range("E26").Value = '  +0.52%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.561'
$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.06'
$ws.Range("E28").Value = '  +0.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.091'
$ws.Range("E29").Value = '  -0.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.14'
$ws.Range("E30").Value = '  +0.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.009'
$ws.Range("E31").Value = '  -0.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09658'
$ws.Range("E32").Value = '  +1.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.547'
$ws.Range("E33").Value = '  -1.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.411'
$ws.Range("E34").Value = '  +3.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.531'
$ws.Range("E35").Value = '  -0.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06583'
$ws.Range("E36").Value = '  +7.75%  '

$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.201'
$ws.Range("E38").Value = '  +1.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5969'
$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.75'
$ws.Range("E40").Value = '  -0.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.944'
$ws.Range("E41").Value = '  -1.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1847'
$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.468'
$ws.Range("E43").Value = '  +3.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.238'
$ws.Range("E44").Value = '  -3.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.34'
$ws.Range("E45").Value = '  -0.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07483'
$ws.Range("E46").Value = '  -1.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5557'
$ws.Range("E47").Value = '  -0.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.984'
$ws.Range("E48").Value = '  +1.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.70'
$ws.Range("E49").Value = '  -0.94%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.52'
$ws.Range("E50").Value = '  +0.45%  '

$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.415'
$ws.Range("E51").Value = '  -0.21%  '
